$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3399.5
$ws.Range("I18").Value = 3399.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3399.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3115.5
$ws.Range("N18").ClearContents()

$ws.Range("H38").Value = 144.6
$ws.Range("I38").Value = 49.555557
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 148.666671
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 223.333329

$ws.Range("H64").Value = 7563.727
$ws.Range("I64").Value = 6520.1
$ws.Range("J64").Value = 18000
$ws.Range("K64").Value = 6520.1
$ws.Range("L64").Value = 18000
$ws.Range("M64").Value = -6272.1
$ws.Range("N64").Value = -18496

$ws.Range("H67").Value = 7563.727
$ws.Range("I67").Value = 6520.1
$ws.Range("J67").Value = 18000
$ws.Range("K67").Value = 6520.1
$ws.Range("L67").Value = 18000
$ws.Range("M67").Value = -5662.1
$ws.Range("N67").Value = -19716

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H138").Value = 7110.6562
$ws.Range("I138").Value = 4878.5
$ws.Range("J138").Value = 8846.777
$ws.Range("K138").Value = 14635.5
$ws.Range("L138").Value = 26540.331
$ws.Range("M138").Value = -9495.5
$ws.Range("N138").Value = -36820.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 458.3
$ws.Range("I2").Value = 476.625
$ws.Range("J2").Value = 385
$ws.Range("K2").Value = 476.625
$ws.Range("L2").Value = 385
$ws.Range("M2").Value = -363.625
$ws.Range("N2").Value = -611

$ws.Range("H11").Value = 349.5
$ws.Range("I11").Value = 299
$ws.Range("J11").Value = 366.33334
$ws.Range("K11").Value = 299
$ws.Range("L11").Value = 366.33334
$ws.Range("M11").Value = -155
$ws.Range("N11").Value = -654.33334

$ws.Range("H32").Value = 4454.143
$ws.Range("I32").Value = 3367.1064
$ws.Range("J32").Value = 29999.5
$ws.Range("K32").Value = 3367.1064
$ws.Range("L32").Value = 29999.5
$ws.Range("M32").Value = -3080.1064

$ws.Range("H116").Value = 458.3
$ws.Range("I116").Value = 476.625
$ws.Range("J116").Value = 385
$ws.Range("K116").Value = 476.625
$ws.Range("L116").Value = 385
$ws.Range("M116").Value = 1817.375
$ws.Range("N116").Value = -4973

$ws.Range("H132").Value = 6832
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 8998
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 26994
$ws.Range("M132").Value = -4970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 458.3
$ws.Range("I3").Value = 476.625
$ws.Range("J3").Value = 385
$ws.Range("K3").Value = 476.625
$ws.Range("L3").Value = 385
$ws.Range("M3").Value = -362.625
$ws.Range("N3").Value = -613

$ws.Range("H105").Value = 4470.8125
$ws.Range("I105").Value = 2065
$ws.Range("J105").Value = 5026
$ws.Range("K105").Value = 2065
$ws.Range("L105").Value = 5026
$ws.Range("M105").Value = -318
$ws.Range("N105").Value = -8520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3048.1428
$ws.Range("I105").Value = 959.25
$ws.Range("J105").Value = 5833.3335
$ws.Range("K105").Value = 959.25
$ws.Range("L105").Value = 5833.3335
$ws.Range("M105").Value = 787.75
$ws.Range("N105").Value = -9327.333500000001

$ws.Range("H107").Value = 657
$ws.Range("I107").Value = 425.16666
$ws.Range("J107").Value = 1352.5
$ws.Range("K107").Value = 425.16666
$ws.Range("L107").Value = 1352.5
$ws.Range("M107").Value = 1494.83334
$ws.Range("N107").Value = -5192.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 120
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 112

$ws.Range("H36").Value = 2706.75
$ws.Range("I36").Value = 173.4
$ws.Range("J36").Value = 6929
$ws.Range("K36").Value = 520.2
$ws.Range("L36").Value = 20787
$ws.Range("M36").Value = -351.2
$ws.Range("N36").Value = -21125

$ws.Range("H59").Value = 1000
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -2460
$ws.Range("N59").ClearContents()

$ws.Range("H139").Value = 4558
$ws.Range("I139").Value = 2800
$ws.Range("J139").Value = 4997.5
$ws.Range("K139").Value = 8400
$ws.Range("L139").Value = 14992.5
$ws.Range("M139").Value = -3260
$ws.Range("N139").Value = -25272.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12515351
$ws.Range("I70").Value = 12515351
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 12515351
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12515081

$ws.Range("H73").Value = 12515351
$ws.Range("I73").Value = 12515351
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 12515351
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -12514415

$ws.Range("H102").Value = 464.30768
$ws.Range("I102").Value = 461.1
$ws.Range("J102").Value = 475
$ws.Range("K102").Value = 461.1
$ws.Range("L102").Value = 475
$ws.Range("M102").Value = 1160.9

$ws.Range("H107").Value = 1056.5454
$ws.Range("I107").Value = 1071.6666
$ws.Range("J107").Value = 1038.4
$ws.Range("K107").Value = 1071.6666
$ws.Range("L107").Value = 1038.4
$ws.Range("M107").Value = 848.3334
$ws.Range("N107").Value = -4878.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 973.1539
$ws.Range("I16").Value = 973.1539
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 973.1539
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -803.1539
$ws.Range("N16").ClearContents()

$ws.Range("H55").Value = 589.25
$ws.Range("I55").Value = 244.81818
$ws.Range("J55").Value = 1010.2222
$ws.Range("K55").Value = 244.81818
$ws.Range("L55").Value = 1010.2222
$ws.Range("M55").Value = -71.81818000000001
$ws.Range("N55").Value = -1356.2222

$ws.Range("H61").Value = 1734.4667
$ws.Range("I61").Value = 1416.8334
$ws.Range("J61").Value = 3005
$ws.Range("K61").Value = 1416.8334
$ws.Range("L61").Value = 3005
$ws.Range("M61").Value = -1214.8334

$ws.Range("H68").Value = 2200
$ws.Range("I68").Value = 2200
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2200
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1451

$ws.Range("H71").Value = 2200
$ws.Range("I71").Value = 2200
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -7256

$ws.Range("H113").Value = 1734.4667
$ws.Range("I113").Value = 1416.8334
$ws.Range("J113").Value = 3005
$ws.Range("K113").Value = 1416.8334
$ws.Range("L113").Value = 3005
$ws.Range("M113").Value = 753.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1067.6
$ws.Range("I81").Value = 1067.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2135.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1074.2
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1067.6
$ws.Range("I84").Value = 1067.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10676
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -5372
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 620.2857
$ws.Range("I113").Value = 390.6
$ws.Range("J113").Value = 1194.5
$ws.Range("K113").Value = 1171.8
$ws.Range("L113").Value = 3583.5
$ws.Range("M113").Value = 998.1999999999998
$ws.Range("N113").Value = -7923.5

$ws.Range("H132").Value = 2265.6956
$ws.Range("I132").Value = 2232.75
$ws.Range("J132").Value = 2485.3333
$ws.Range("K132").Value = 6698.25
$ws.Range("L132").Value = 7455.999899999999
$ws.Range("M132").Value = -4168.25

$ws.Range("H136").Value = 1738.1111
$ws.Range("I136").Value = 1235.9286
$ws.Range("J136").Value = 3495.75
$ws.Range("K136").Value = 3707.7858
$ws.Range("L136").Value = 10487.25
$ws.Range("M136").Value = -1157.7858
$ws.Range("N136").Value = -15587.25
